$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test-Cases")
$ws.Activate()

# Add a new column (J) carrying "Test Sep Value" alongside the existing
# "Rejected" values in column I, for the two populated test-case rows (6 & 12)
$ws.Range("J6").Value = "Test Sep Value"
$ws.Range("J12").Value = "Test Sep Value"

# Update the view to match the saved state: scrolled right a couple of
# columns (was F1, now H1) with the active cell moved from I12 to J12
$ws.Range("J12").Select()
$excel.ActiveWindow.ScrollColumn = 8
